$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("doSignIn")

# --- Hyperlinks -------------------------------------------------------
# The sheet's hyperlink collection in this engine can only be wiped in
# one shot (deleting via any Range.Hyperlinks.Delete() clears the whole
# worksheet) and rebuilt with Hyperlinks.Add — so remove everything and
# recreate the links that should remain (the old B3 link is dropped),
# in the exact order that reproduces the desired rId numbering.
$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:QqwertyQ@123!", "", "", "QqwertyQ@123!") | Out-Null
$ws.Range("B2").Style = "Hyperlink"
$ws.Range("B2").NumberFormat = "@"

$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:abctestemail2372!!!!@gmail.com") | Out-Null
$ws.Range("A6").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("A2:A5"), "mailto:abctestemail237!!@gmail.com", "", "", "abctestemail237!!@gmail.com") | Out-Null
$ws.Range("A2:A5").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:abc1@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:abctestemail2379!!!!@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:abctestemail2370!!!!@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:abctestemail2371!!!!@gmail.com") | Out-Null
$ws.Range("A2:A6").Style = "Hyperlink"

# Re-adding the hyperlinks above can stomp the displayed text of the
# cells that carry an explicit display string (B2, A2) — put the
# untouched values (row 2 does not change in this edit) back exactly.
$ws.Range("B2").Value = "AydenLiam1213"
$ws.Range("A2").Value = "abc1@gmail.com"

# --- New sign-in test data (rows 3-6) ---------------------------------
# Written in the same order the new strings were originally appended to
# the shared-string table: A6, A3, A5, then A4.
$ws.Range("A6").Value = "abctestemail4!!!!@gmail.com"
$ws.Range("A3").Value = 'abctestemail1$@gmail.com'
$ws.Range("A5").Value = 'abctestemail3$@gmail.com'
$ws.Range("A4").Value = "james.corley@gmail.com"

$ws.Range("B3").Value = "AydenLiam1213"
$ws.Range("B4").Value = "AydenLiam1213"
$ws.Range("B5").Value = "AydenLiam1213"
$ws.Range("B6").Value = "AydenLiam1213"

# --- Selection ----------------------------------------------------------
$ws.Range("A4").Select() | Out-Null
